# The project this template ships with was renamed (TEDU/TeduCore -> TanoApp),
# so the first worksheet ("TEDUOrder") is renamed to "TanoOrder" to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "TanoOrder"
